# GuildName.xlsx edit: unify the conception of DataNode, DataTable, Entity.
# The only content-level change in this revision is renaming the sheet
# that used to be called "Property" to "DataNode", and the cursor/
# selection position that was left on the sheet when it was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename "Property" -> "DataNode"
$ws.Name = "DataNode"

# Leave the saved selection where the author left it (D39) instead of
# the old A9 default.
$ws.Range("D39").Select() | Out-Null
